$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Updated query timestamps for the "time_taken" column (F) on the data sheet,
# reflecting a re-run of the panel query.
$newTimes = @(
    "2021-10-05 14:35:47.601140",
    "2021-10-05 14:35:47.601147",
    "2021-10-05 14:35:47.601150",
    "2021-10-05 14:35:47.601153",
    "2021-10-05 14:35:47.601156",
    "2021-10-05 14:35:47.601159",
    "2021-10-05 14:35:47.601161",
    "2021-10-05 14:35:47.601164",
    "2021-10-05 14:35:47.601167",
    "2021-10-05 14:35:47.601169",
    "2021-10-05 14:35:47.601171",
    "2021-10-05 14:35:47.601174",
    "2021-10-05 14:35:47.601176",
    "2021-10-05 14:35:47.601179",
    "2021-10-05 14:35:47.601181",
    "2021-10-05 14:35:47.601184",
    "2021-10-05 14:35:47.601186",
    "2021-10-05 14:35:47.601189",
    "2021-10-05 14:35:47.601192",
    "2021-10-05 14:35:47.601194",
    "2021-10-05 14:35:47.601196",
    "2021-10-05 14:35:47.601199",
    "2021-10-05 14:35:47.601201",
    "2021-10-05 14:35:47.601204",
    "2021-10-05 14:35:47.601206",
    "2021-10-05 14:35:47.601209",
    "2021-10-05 14:35:47.601211",
    "2021-10-05 14:35:47.601214",
    "2021-10-05 14:35:47.601217",
    "2021-10-05 14:35:47.601219",
    "2021-10-05 14:35:47.601222",
    "2021-10-05 14:35:47.601224",
    "2021-10-05 14:35:47.601227",
    "2021-10-05 14:35:47.601229",
    "2021-10-05 14:35:47.601232",
    "2021-10-05 14:35:47.601234",
    "2021-10-05 14:35:47.601237",
    "2021-10-05 14:35:47.601239",
    "2021-10-05 14:35:47.601242",
    "2021-10-05 14:35:47.601244",
    "2021-10-05 14:35:47.601247",
    "2021-10-05 14:35:47.601250",
    "2021-10-05 14:35:47.601252",
    "2021-10-05 14:35:47.601255",
    "2021-10-05 14:35:47.601257",
    "2021-10-05 14:35:47.601260",
    "2021-10-05 14:35:47.601262",
    "2021-10-05 14:35:47.601265",
    "2021-10-05 14:35:47.601267",
    "2021-10-05 14:35:47.601270",
    "2021-10-05 14:35:47.601272",
    "2021-10-05 14:35:47.601274"
)

for ($i = 0; $i -lt $newTimes.Count; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# Add a new "metadata" worksheet, placed after "data", describing the panel query.
$ws2 = $wb.Worksheets.Add($null, $dataSheet)
$ws2.Name = "metadata"

$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Stroke"
$ws2.Range("C2").Value = 3141
$ws2.Range("E2").Value = "2021-08-11T07:39:17.004134Z"
$ws2.Range("F2").Value = "2021-10-05 14:35:47.597434"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3141/?format=json"

# data_version must be stored as text "1.6" (not the number 1.6).
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "1.6"
$ws2.Range("D2").Style = $ws2.Range("B2").Style

# Match the header row (B1:G1) and the index cell (A2) to the bold/bordered
# header style used on the "data" sheet.
$dataSheet.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

$dataSheet.Select()
